$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-sort the "Periodo Mora" (col E) rows into ascending order (1611..1705),
# keeping each period's "Valor Mora" (col F) attached, and refresh the
# "Salario Basico" (col G) to the new value for all periods.
$ws.Range("E16").Value = "1611"
$ws.Range("F16").Value = 27578
$ws.Range("G16").Value = 781242

$ws.Range("E17").Value = "1612"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 781242

$ws.Range("E18").Value = "1701"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 781242

$ws.Range("E19").Value = "1704"
$ws.Range("F19").Value = 29509
$ws.Range("G19").Value = 781242

$ws.Range("E20").Value = "1705"
$ws.Range("F20").Value = 29509
$ws.Range("G20").Value = 781242
